# Update cryptos list: prices and 1h-volume percentages refreshed,
# plus three name/link pairs that swapped rank order (rows 16/17, 32/33, 45/46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.636.58"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "3.213.70"
$ws.Range("E3").Value = "  +1.63%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "218.64"
$ws.Range("E5").Value = "  +6.39%  "

$ws.Range("D6").Value = "656.13"
$ws.Range("E6").Value = "  +7.90%  "

$ws.Range("D7").Value = "0.399"
$ws.Range("E7").Value = "  +5.70%  "

$ws.Range("D8").Value = "0.698"
$ws.Range("E8").Value = "  +5.50%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "3.212.90"
$ws.Range("E10").Value = "  +1.80%  "

$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +8.53%  "

$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +5.81%  "

$ws.Range("D14").Value = "5.43"
$ws.Range("E14").Value = "  +3.46%  "

$ws.Range("D15").Value = "33.75"
$ws.Range("E15").Value = "  +5.08%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "90.212.39"
$ws.Range("E16").Value = "  +3.54%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.796.83"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").Value = "3.221.38"
$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("D19").Value = "3.42"
$ws.Range("E19").Value = "  +14.57%  "

$ws.Range("D20").Value = "0.0000229"
$ws.Range("E20").Value = "  +76.12%  "

$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  +1.78%  "

$ws.Range("D22").Value = "441.78"
$ws.Range("E22").Value = "  +6.78%  "

$ws.Range("D23").Value = "8.72"
$ws.Range("E23").Value = "  +2.85%  "

$ws.Range("D24").Value = "5.13"
$ws.Range("E24").Value = "  +0.87%  "

$ws.Range("D25").Value = "5.36"
$ws.Range("E25").Value = "  +3.81%  "

$ws.Range("D26").Value = "12.04"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").Value = "82.17"
$ws.Range("E27").Value = "  +12.13%  "

$ws.Range("D28").Value = "3.372.18"
$ws.Range("E28").Value = "  +1.54%  "

$ws.Range("D30").Value = "0.161"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "550.92"
$ws.Range("E32").Value = "  +2.22%  "

$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "4.07"
$ws.Range("E33").Value = "  +35.87%  "

$ws.Range("D34").Value = "8.55"

$ws.Range("D35").Value = "7.13"
$ws.Range("E35").Value = "  +6.38%  "

$ws.Range("D36").Value = "1.95"
$ws.Range("E36").Value = "  +6.39%  "

$ws.Range("E37").Value = "  +1.07%  "

$ws.Range("D38").Value = "22.65"
$ws.Range("E38").Value = "  +3.73%  "

$ws.Range("D39").Value = "22.42"
$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("E40").Value = "  -3.02%  "

$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("E42").Value = "  +3.23%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("E44").Value = "  +1.73%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "45.91"
$ws.Range("E45").Value = "  +6.14%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "146.63"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").Value = "174.77"
$ws.Range("E47").Value = "  +1.30%  "

$ws.Range("D48").Value = "0.768"
$ws.Range("E48").Value = "  +10.57%  "

$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("D50").Value = "1.26"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("D51").Value = "0.627"
$ws.Range("E51").Value = "  +6.95%  "

